$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.685.55'
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('D3').Value = '3.797.25'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'697.58"
$ws.Range('E5').Value = '  +8.04%  '
$ws.Range('D6').Value = "'172.47"
$ws.Range('E6').Value = '  +3.88%  '
$ws.Range('D7').Value = '3.795.44'
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').Value = "'0.162"
$ws.Range('E10').Value = '  +1.91%  '
$ws.Range('D11').Value = "'7.25"
$ws.Range('E11').Value = '  +5.14%  '
$ws.Range('D12').Value = "'0.458"
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('D13').Value = "'0.0000258"
$ws.Range('E13').Value = '  +7.58%  '
$ws.Range('D14').Value = "'36.14"
$ws.Range('D15').Value = '4.443.79'
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('D16').Value = '3.786.13'
$ws.Range('E16').Value = '  +0.96%  '
$ws.Range('D17').Value = '70.717.34'
$ws.Range('E17').Value = '  +1.70%  '
$ws.Range('D18').Value = "'17.73"
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').Value = "'7.17"
$ws.Range('E19').Value = '  +2.01%  '
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('D21').Value = "'11.16"
$ws.Range('E21').Value = '  +16.20%  '
$ws.Range('D22').Value = "'478.03"
$ws.Range('E22').Value = '  +1.94%  '
$ws.Range('D23').Value = "'0.709"
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').Value = "'83.76"
$ws.Range('E24').Value = '  +2.06%  '
$ws.Range('D25').Value = "'0.0000142"
$ws.Range('E25').Value = '  -1.13%  '
$ws.Range('D26').Value = "'12.30"
$ws.Range('E26').Value = '  -0.33%  '
$ws.Range('D27').Value = "'2.15"
$ws.Range('E27').Value = '  +1.88%  '
$ws.Range('D28').Value = "'10.39"
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('D29').Value = '3.949.81'
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('E30').Value = '  -0.14%  '
$ws.Range('D31').Value = "'3.11"
$ws.Range('E31').Value = '  +15.00%  '
$ws.Range('E32').Value = '  +0.47%  '
$ws.Range('E33').Value = '  +4.55%  '
$ws.Range('D34').Value = "'0.189"
$ws.Range('E34').Value = '  +9.06%  '
$ws.Range('D35').Value = "'29.37"
$ws.Range('E35').Value = '  +1.77%  '
$ws.Range('D36').Value = "'9.22"
$ws.Range('E36').Value = '  +3.78%  '
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('E38').Value = '  +2.27%  '
$ws.Range('D39').Value = "'3.40"
$ws.Range('E39').Value = '  +3.23%  '
$ws.Range('D40').Value = "'5.99"
$ws.Range('E40').Value = '  +2.24%  '
$ws.Range('D41').Value = "'2.25"
$ws.Range('E41').Value = '  +12.66%  '
$ws.Range('D42').Value = "'0.974"
$ws.Range('E42').Value = '  +1.76%  '
$ws.Range('E43').Value = '  +20.45%  '
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = "'163.74"
$ws.Range('E46').Value = '  +4.10%  '
$ws.Range('D47').Value = "'48.85"
$ws.Range('E47').Value = '  +2.34%  '
$ws.Range('D48').Value = "'44.31"
$ws.Range('E48').Value = '  -2.43%  '
$ws.Range('D49').Value = "'0.299"
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').Value = "'1.37"
$ws.Range('E50').Value = '  -1.69%  '
$ws.Range('D51').Value = "'8.56"
$ws.Range('E51').Value = '  +2.08%  '
